$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bug")

# Copy formatting from the last existing data row (row 15) down to the new row 16
$ws.Range("A15:D15").Copy() | Out-Null
$ws.Range("A16:D16").PasteSpecial(-4122) | Out-Null

# Populate the new row (Bug/Task No 15: a new Code task)
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "Correction of All masters page and complete CURD operation"
$ws.Range("C16").Value = "Code"
$ws.Range("D16").Value = "In-Progress"

# Update status of rows 10 and 11 (currently "Pending") to "In-Progress"
$ws.Range("D10").Value = "In-Progress"
$ws.Range("D11").Value = "In-Progress"

# Update the active selection to match the edited cell
$ws.Range("D10").Select()
